$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the extra row (previously row 4 / "125M") - only 2 data rows remain
$ws.Rows.Item(4).Delete()

# Row 1 - headers
$ws.Range("A1").Value = "model"
$ws.Range("B1").Value = "system"
$ws.Range("C1").Value = "category"
$ws.Range("D1").Value = "symptomTitle"
$ws.Range("E1").Value = "diagnosisTreeId"
$ws.Range("F1").Value = "diagnosisResultId"
$ws.Range("G1").Value = "title"
$ws.Range("H1").Value = "description"
$ws.Range("I1").Value = "fixSteps"
$ws.Range("J1").Value = "tags"
$ws.Range("K1").Value = "references"
$ws.Range("L1").Value = "parts"
$ws.Range("M1").Value = "photo_1"
$ws.Range("N1").Value = "photo_1_desc"
$ws.Range("O1").Value = "photo_2"
$ws.Range("P1").Value = "photo_2_desc"
$ws.Range("Q1").Value = "photo_3"
$ws.Range("R1").Value = "photo_3_desc"
$ws.Range("S1").Value = "photo_4"
$ws.Range("T1").Value = "photo_4_desc"
$ws.Range("U1").Value = "photo_5"
$ws.Range("V1").Value = "photo_5_desc"

# Row 2 - 350D / engine case (A2 "350D" unchanged - leave as-is)
$ws.Range("B2").Value = "engine"
$ws.Range("C2").Value = "Engine"
$ws.Range("D2").Value = "Poor acceleration"
$ws.Range("E2").Value = "poor_acceleration_v1"
$ws.Range("F2").Value = "r3"
$ws.Range("G2").Value = "Throttle cable free play out of spec"
$ws.Range("H2").Value = "Low acceleration after warm-up"
$ws.Range("I2").Value = "Adjust throttle cable free play to spec"
$ws.Range("J2").Value = "engine,acceleration"
$ws.Range("K2").Value = "https://example.com"
$ws.Range("L2").Value = "Throttle cable"
# M2, N2 were already blank and remain blank - leave as-is
# O2:V2 are brand-new columns for this row - leave blank (no prior content)

# Row 3 - 368G / engine idle case (A3 "368G" unchanged - leave as-is)
$ws.Range("B3").Value = "engine"
$ws.Range("C3").Value = "Idle"
$ws.Range("D3").Value = "High hot idle"
$ws.Range("E3").Value = "high_hot_idle_v1"
$ws.Range("F3").Value = "r4"
$ws.Range("G3").Value = "Coolant temp sensor fault"
$ws.Range("H3").Value = "Idle speed remains high after warm-up"
$ws.Range("I3").Value = "Inspect wiring and replace coolant temperature sensor"
$ws.Range("J3").Value = "idle,sensor"
# K3 was already blank and remains blank - leave as-is
$ws.Range("L3").Value = "Coolant temperature sensor"
# M3, N3 were already blank and remain blank - leave as-is
# O3:V3 are brand-new columns for this row - leave blank (no prior content)
